# Apply updated TPM-derived statistics to the Bmp2-Eng LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.729584333333333
$ws.Range("H2").Value = 5.188753
$ws.Range("I2").Value = 0.2476387648475193
$ws.Range("J2").Value = 0.2476387648475193
$ws.Range("M2").Value = 247.0944516666667
$ws.Range("N2").Value = 741.283355
$ws.Range("O2").Value = 0.8050739182622993
$ws.Range("P2").Value = 0.8050739182622993
$ws.Range("Q2").Value = 427.3706924562572
$ws.Range("R2").Value = 3846.336232106315
$ws.Range("S2").Value = 0.1993675107294285
$ws.Range("T2").Value = 0.1993675107294285

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.729584333333333
$ws.Range("H3").Value = 5.188753
$ws.Range("I3").Value = 0.2476387648475193
$ws.Range("J3").Value = 0.2476387648475193
$ws.Range("O3").Value = 0.1379009747488701
$ws.Range("P3").Value = 0.13790097474887
$ws.Range("Q3").Value = 73.20425333865556
$ws.Range("R3").Value = 658.8382800479001
$ws.Range("S3").Value = 0.03414962705807913
$ws.Range("T3").Value = 0.03414962705807912

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.729584333333333
$ws.Range("H4").Value = 5.188753
$ws.Range("I4").Value = 0.2476387648475193
$ws.Range("J4").Value = 0.2476387648475193
$ws.Range("M4").Value = 11.590146
$ws.Range("N4").Value = 34.770438
$ws.Range("O4").Value = 0.03776258103132013
$ws.Range("P4").Value = 0.03776258103132013
$ws.Range("Q4").Value = 20.046134942646
$ws.Range("R4").Value = 180.415214483814
$ws.Range("S4").Value = 0.009351478924050478
$ws.Range("T4").Value = 0.009351478924050478

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.729584333333333
$ws.Range("H5").Value = 5.188753
$ws.Range("I5").Value = 0.2476387648475193
$ws.Range("J5").Value = 0.2476387648475193
$ws.Range("M5").Value = 5.912082333333333
$ws.Range("N5").Value = 17.736247
$ws.Range("O5").Value = 0.01926252595751047
$ws.Range("P5").Value = 0.01926252595751047
$ws.Range("Q5").Value = 10.22544498111011
$ws.Range("R5").Value = 92.029004829991
$ws.Range("S5").Value = 0.004770148135961172
$ws.Range("T5").Value = 0.004770148135961172

# Row 6
$ws.Range("I6").Value = 0.2307941364328804
$ws.Range("J6").Value = 0.2307941364328804
$ws.Range("M6").Value = 247.0944516666667
$ws.Range("N6").Value = 741.283355
$ws.Range("O6").Value = 0.8050739182622993
$ws.Range("P6").Value = 0.8050739182622993
$ws.Range("Q6").Value = 398.3005244065773
$ws.Range("R6").Value = 3584.704719659195
$ws.Range("S6").Value = 0.1858063397299827
$ws.Range("T6").Value = 0.1858063397299827

# Row 7
$ws.Range("I7").Value = 0.2307941364328804
$ws.Range("J7").Value = 0.2307941364328804
$ws.Range("O7").Value = 0.1379009747488701
$ws.Range("P7").Value = 0.13790097474887
$ws.Range("S7").Value = 0.03182673638041791
$ws.Range("T7").Value = 0.0318267363804179

# Row 8
$ws.Range("I8").Value = 0.2307941364328804
$ws.Range("J8").Value = 0.2307941364328804
$ws.Range("M8").Value = 11.590146
$ws.Range("N8").Value = 34.770438
$ws.Range("O8").Value = 0.03776258103132013
$ws.Range("P8").Value = 0.03776258103132013
$ws.Range("Q8").Value = 18.682577446038
$ws.Range("R8").Value = 168.143197014342
$ws.Range("S8").Value = 0.008715382278600199
$ws.Range("T8").Value = 0.008715382278600199

# Row 9
$ws.Range("I9").Value = 0.2307941364328804
$ws.Range("J9").Value = 0.2307941364328804
$ws.Range("M9").Value = 5.912082333333333
$ws.Range("N9").Value = 17.736247
$ws.Range("O9").Value = 0.01926252595751047
$ws.Range("P9").Value = 0.01926252595751047
$ws.Range("Q9").Value = 9.529900318758111
$ws.Range("R9").Value = 85.769102868823
$ws.Range("S9").Value = 0.004445678043879572
$ws.Range("T9").Value = 0.004445678043879572

# Row 10
$ws.Range("G10").Value = 2.743651333333334
$ws.Range("H10").Value = 8.230954000000001
$ws.Range("I10").Value = 0.3928310486309039
$ws.Range("J10").Value = 0.3928310486309038
$ws.Range("M10").Value = 247.0944516666667
$ws.Range("N10").Value = 741.283355
$ws.Range("O10").Value = 0.8050739182622993
$ws.Range("P10").Value = 0.8050739182622993
$ws.Range("Q10").Value = 677.9410217745191
$ws.Range("R10").Value = 6101.469195970671
$ws.Range("S10").Value = 0.3162580315363697
$ws.Range("T10").Value = 0.3162580315363696

# Row 11
$ws.Range("G11").Value = 2.743651333333334
$ws.Range("H11").Value = 8.230954000000001
$ws.Range("I11").Value = 0.3928310486309039
$ws.Range("J11").Value = 0.3928310486309038
$ws.Range("O11").Value = 0.1379009747488701
$ws.Range("P11").Value = 0.13790097474887
$ws.Range("Q11").Value = 116.1244024980222
$ws.Range("R11").Value = 1045.1196224822
$ws.Range("S11").Value = 0.05417178451782242
$ws.Range("T11").Value = 0.0541717845178224

# Row 12
$ws.Range("G12").Value = 2.743651333333334
$ws.Range("H12").Value = 8.230954000000001
$ws.Range("I12").Value = 0.3928310486309039
$ws.Range("J12").Value = 0.3928310486309038
$ws.Range("M12").Value = 11.590146
$ws.Range("N12").Value = 34.770438
$ws.Range("O12").Value = 0.03776258103132013
$ws.Range("P12").Value = 0.03776258103132013
$ws.Range("Q12").Value = 31.799319526428
$ws.Range("R12").Value = 286.193875737852
$ws.Range("S12").Value = 0.01483431430554297
$ws.Range("T12").Value = 0.01483431430554297

# Row 13
$ws.Range("G13").Value = 2.743651333333334
$ws.Range("H13").Value = 8.230954000000001
$ws.Range("I13").Value = 0.3928310486309039
$ws.Range("J13").Value = 0.3928310486309038
$ws.Range("M13").Value = 5.912082333333333
$ws.Range("N13").Value = 17.736247
$ws.Range("O13").Value = 0.01926252595751047
$ws.Range("P13").Value = 0.01926252595751047
$ws.Range("Q13").Value = 16.22069257662645
$ws.Range("R13").Value = 145.986233189638
$ws.Range("S13").Value = 0.007566918271168845
$ws.Range("T13").Value = 0.007566918271168844

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.8991316666666668
$ws.Range("H14").Value = 2.697395
$ws.Range("I14").Value = 0.1287360500886965
$ws.Range("J14").Value = 0.1287360500886965
$ws.Range("M14").Value = 247.0944516666667
$ws.Range("N14").Value = 741.283355
$ws.Range("O14").Value = 0.8050739182622993
$ws.Range("P14").Value = 0.8050739182622993
$ws.Range("Q14").Value = 222.1704461511362
$ws.Range("R14").Value = 1999.534015360225
$ws.Range("S14").Value = 0.1036420362665185
$ws.Range("T14").Value = 0.1036420362665185

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8991316666666668
$ws.Range("H15").Value = 2.697395
$ws.Range("I15").Value = 0.1287360500886965
$ws.Range("J15").Value = 0.1287360500886965
$ws.Range("O15").Value = 0.1379009747488701
$ws.Range("P15").Value = 0.13790097474887
$ws.Range("Q15").Value = 38.05553799427778
$ws.Range("R15").Value = 342.4998419485
$ws.Range("S15").Value = 0.01775282679255061
$ws.Range("T15").Value = 0.0177528267925506

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.8991316666666668
$ws.Range("H16").Value = 2.697395
$ws.Range("I16").Value = 0.1287360500886965
$ws.Range("J16").Value = 0.1287360500886965
$ws.Range("M16").Value = 11.590146
$ws.Range("N16").Value = 34.770438
$ws.Range("O16").Value = 0.03776258103132013
$ws.Range("P16").Value = 0.03776258103132013
$ws.Range("Q16").Value = 10.42106728989
$ws.Range("R16").Value = 93.78960560901001
$ws.Range("S16").Value = 0.004861405523126489
$ws.Range("T16").Value = 0.004861405523126489

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8991316666666668
$ws.Range("H17").Value = 2.697395
$ws.Range("I17").Value = 0.1287360500886965
$ws.Range("J17").Value = 0.1287360500886965
$ws.Range("M17").Value = 5.912082333333333
$ws.Range("N17").Value = 17.736247
$ws.Range("O17").Value = 0.01926252595751047
$ws.Range("P17").Value = 0.01926252595751047
$ws.Range("Q17").Value = 5.315740441840556
$ws.Range("R17").Value = 47.841663976565
$ws.Range("S17").Value = 0.002479781506500885
$ws.Range("T17").Value = 0.002479781506500885

